# EPBDS-7507 Workaround for Shared Formulas in xls files.
# Adds a new "Formulas" sheet after the existing sheets, populates B3:B12
# with numbers 1..10 and C3:C12 with a shared formula "=B*10", then
# selects C13 on that sheet and makes it the active tab.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Formulas"

# Fill B3:B12 with 1..10.
for ($i = 0; $i -lt 10; $i++) {
    $row = 3 + $i
    $newSheet.Range("B$row").Value = $i + 1
}

# C3:C12 = B*10, written as one shared formula across the range.
$newSheet.Range("C3:C12").Formula = "=B3*10"

# Activate the new sheet and place the selection on C13, matching the
# authored workbook's view state.
$newSheet.Activate()
$newSheet.Range("C13").Select()
